# Add a new results row (row 74) to the "Results" sheet, duplicating the
# previous run's row (row 73) but with a new timestamp. This mirrors a
# re-run of the MARS m3c2 stats pipeline (see commit message: the "mov"
# point-cloud is now included in the statistics computation) that produced
# an additional logged entry in the distances/statistics table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$row = 74

# A: Timestamp, B: Folder, C: Version (text columns)
$ws.Cells.Item($row, 1).Value  = '2025-09-02 11:39:20'
$ws.Cells.Item($row, 2).Value  = '0342-0349'
$ws.Cells.Item($row, 3).Value  = 'mov-ref'

# D: Total Points .. AG: Outlier Threshold (numeric columns)
$ws.Cells.Item($row, 4).Value  = 709128
$ws.Cells.Item($row, 5).Value  = 0.0622794074939492
$ws.Cells.Item($row, 6).Value  = 0.124558814987898
$ws.Cells.Item($row, 7).Value  = 9634
$ws.Cells.Item($row, 8).Value  = 0.01358569961981476
$ws.Cells.Item($row, 9).Value  = 0.9864143003801853
$ws.Cells.Item($row, 10).Value = 699494
$ws.Cells.Item($row, 11).Value = -621.7036519999999
$ws.Cells.Item($row, 12).Value = 143.812791320498
$ws.Cells.Item($row, 13).Value = 683858
$ws.Cells.Item($row, 14).Value = -1245.446727
$ws.Cells.Item($row, 15).Value = 89.512016083193
$ws.Cells.Item($row, 16).Value = -0.11556
$ws.Cells.Item($row, 17).Value = 0.132246
$ws.Cells.Item($row, 18).Value = -0.0008887905428781375
$ws.Cells.Item($row, 19).Value = -0.002101
$ws.Cells.Item($row, 20).Value = 0.01433860038117433
$ws.Cells.Item($row, 21).Value = 0.01431102764520784
$ws.Cells.Item($row, 22).Value = 0.009318774988777602
$ws.Cells.Item($row, 23).Value = 0.008409307199999999
$ws.Cells.Item($row, 24).Value = -0.043002
$ws.Cells.Item($row, 25).Value = 0.043009
$ws.Cells.Item($row, 26).Value = -0.001821206635003173
$ws.Cells.Item($row, 27).Value = -0.002249
$ws.Cells.Item($row, 28).Value = 0.01144083450710502
$ws.Cells.Item($row, 29).Value = 0.0112949502261668
$ws.Cells.Item($row, 30).Value = 0.008215368975723031
$ws.Cells.Item($row, 31).Value = 0.008176539
$ws.Cells.Item($row, 32).Value = 3
$ws.Cells.Item($row, 33).Value = 0.04301580114352299

# AH: Outlier Method (text column)
$ws.Cells.Item($row, 34).Value = 'rmse'

# AI: Inlier Count .. BK: Kurtosis (numeric columns)
$ws.Cells.Item($row, 35).Value = 683858
$ws.Cells.Item($row, 36).Value = 264076
$ws.Cells.Item($row, 37).Value = 419712
$ws.Cells.Item($row, 38).Value = 12902
$ws.Cells.Item($row, 39).Value = 2734
$ws.Cells.Item($row, 40).Value = 15636
$ws.Cells.Item($row, 41).Value = 0.03989147320286518
$ws.Cells.Item($row, 42).Value = 0.04337597406792291
$ws.Cells.Item($row, 43).Value = -0.018882
$ws.Cells.Item($row, 44).Value = -0.007848000000000001
$ws.Cells.Item($row, 45).Value = 0.003504
$ws.Cells.Item($row, 46).Value = 0.024597
$ws.Cells.Item($row, 47).Value = 0.011352
$ws.Cells.Item($row, 48).Value = -0.018424
$ws.Cells.Item($row, 49).Value = -0.007889
$ws.Cells.Item($row, 50).Value = 0.003141
$ws.Cells.Item($row, 51).Value = 0.018931
$ws.Cells.Item($row, 52).Value = 0.01103
$ws.Cells.Item($row, 53).Value = -0.0008887905428781375
$ws.Cells.Item($row, 54).Value = 0.01431102764520784
$ws.Cells.Item($row, 55).Value = 238689480875.462
$ws.Cells.Item($row, 56).Value = 5.675151041459209
$ws.Cells.Item($row, 57).Value = 0.1035729672089861
$ws.Cells.Item($row, 58).Value = -0.09813477944421571
$ws.Cells.Item($row, 59).Value = 0.001960367063261137
$ws.Cells.Item($row, 60).Value = -0.338411737141495
$ws.Cells.Item($row, 61).Value = 9860681782068.098
$ws.Cells.Item($row, 62).Value = 1.553251175032938
$ws.Cells.Item($row, 63).Value = 7.219025590519583

# BL: Distances Path, BM: Params Path (text columns)
$ws.Cells.Item($row, 64).Value = 'data\0342-0349\python_mov-ref_m3c2_distances.txt'
$ws.Cells.Item($row, 65).Value = 'data\0342-0349\python_mov-ref_m3c2_params.txt'
